$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.459.11"
$ws.Range("E2").Value = "  +0.54%  "

$ws.Range("D3").Value = "1.875.95"
$ws.Range("E3").Value = "  -0.42%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.021"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.48"
$ws.Range("E5").Value = "  +0.28%  "

$ws.Range("E6").Value = "  +0.58%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5118"
$ws.Range("E7").Value = "  -0.28%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3961"
$ws.Range("E8").Value = "  +0.93%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08446"
$ws.Range("E9").Value = "  +1.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.111"
$ws.Range("E10").Value = "  -1.02%  "

$ws.Range("B11").Value = "Polkadot"
$ws.Range("C11").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.248"
$ws.Range("E11").Value = "  +0.24%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.875.38"
$ws.Range("E12").Value = "  -0.74%  "

$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.47"
$ws.Range("E13").Value = "  -0.22%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.241"
$ws.Range("E14").Value = "  -0.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.021"
$ws.Range("E15").Value = "  +0.69%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001111"
$ws.Range("E16").Value = "  +0.58%  "

$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.01"
$ws.Range("E17").Value = "  -0.08%  "

$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06775"
$ws.Range("E18").Value = "  +1.02%  "

$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.73"
$ws.Range("E19").Value = "  -0.39%  "

$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.020"
$ws.Range("E20").Value = "  +0.72%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.943"
$ws.Range("E21").Value = "  -1.32%  "

$ws.Range("B22").Value = "WrappedBTC"
$ws.Range("C22").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D22").Value = "28.498.58"
$ws.Range("E22").Value = "  +0.53%  "

$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.16"
$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.293"
$ws.Range("E24").Value = "  +0.80%  "

$ws.Range("B25").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C25").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D25").Value = "2.088.77"
$ws.Range("E25").Value = "  -0.68%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.67"
$ws.Range("E26").Value = "  +0.64%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.75"
$ws.Range("E27").Value = "  -0.18%  "

$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.349"
$ws.Range("E28").Value = "  -3.86%  "

$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "127.01"
$ws.Range("E29").Value = "  +0.29%  "

$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.1055"
$ws.Range("E30").Value = "  -0.40%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.039"
$ws.Range("E31").Value = "  +0.21%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.760"
$ws.Range("E32").Value = "  -1.71%  "

$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.643"
$ws.Range("E33").Value = "  +0.47%  "

$ws.Range("B34").Value = "VeChain"
$ws.Range("C34").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.02435"
$ws.Range("E34").Value = "  -0.38%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06469"
$ws.Range("E35").Value = "  -1.35%  "

$ws.Range("B36").Value = "Algorand"
$ws.Range("C36").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2177"
$ws.Range("E36").Value = "  -1.45%  "

$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.819"
$ws.Range("E37").Value = "  -6.72%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.264"
$ws.Range("E38").Value = "  +1.89%  "

$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.184"
$ws.Range("E39").Value = "  -0.29%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6385"
$ws.Range("E40").Value = "  -1.54%  "

$ws.Range("B41").Value = "InternetComputer(DFINITY)"
$ws.Range("C41").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.988"
$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.22"
$ws.Range("E42").Value = "  +0.14%  "

$ws.Range("B43").Value = "Decentraland"
$ws.Range("C43").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6050"
$ws.Range("E43").Value = "  -0.74%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.02"
$ws.Range("E44").Value = "  -0.30%  "

$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.714"
$ws.Range("E45").Value = "  +0.29%  "

$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.992"
$ws.Range("E46").Value = "  -1.07%  "

$ws.Range("B47").Value = "WEMIXTOKEN"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.205"
$ws.Range("E47").Value = "  -6.39%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.25"
$ws.Range("E48").Value = "  +0.87%  "

$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.205"
$ws.Range("E49").Value = "  -2.65%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06842"
$ws.Range("E50").Value = "  -1.02%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "76.57"
$ws.Range("E51").Value = "  -1.75%  "

